$d = $word.ActiveDocument

# Locate the target paragraph: the empty ("_GoBack" bookmark) paragraph that
# currently carries ListParagraph/numPr (ilvl 0, numId 18) plus bold/italic,
# immediately after "Stel voor of het handig is om te bouwen".
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $prevText = ""
    if ($i -gt 1) { $prevText = $d.Paragraphs.Item($i - 1).Range.Text }
    if ($para.Range.Text -eq "`r" -and $para.Range.Bold -and $para.Range.Italic -and $prevText -like "*Stel voor of het handig is om te bouwen*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:i/><w:lang w:val="nl-NL"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="nl-NL"/></w:rPr><w:t>Docentmeeting 20 november 2017</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Status: documentatie is geupdate (SRS en SRD), deze week begonnen met code schrijven, opgesplitst in twee back-end en twee front-ends, repository, Play Framework, zit een testmodule in verwerkt, voor front-end AngularJS</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Testingtool: Managementtools </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> GitHub issues aan toevoegen, prioriseren etc...managementtooling </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> We maken gebruik van Trello</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Definition of done: </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>definiëeren binnen een sprint van 1-2 weken, dat iedere iteratie een product levert</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Back-end: REST service</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:sym w:font="Wingdings" w:char="F0DF"/></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> unit test</w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> --</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">HTTP-request </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:sym w:font="Wingdings" w:char="F0DF"/></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Integration test</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Front-end: AngularJS </w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:sym w:font="Wingdings" w:char="F0DF"/></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> UI getest worden </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> gebruik maken van</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:lang w:val="nl-NL"/></w:rPr><w:t>Selendium</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Iteration test </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>en</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> unit test op features, op de business logic</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Features implementeren binnen 1-2 weken, aan de hand van requirementstesting / integration testing / unitttesting </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> functionaliteit als done geevalueerd kan worden.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Volgende week Product Owner rol</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Presenteren testing en waarom welke keuzes</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Implementaties</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Huidige features die we implementeren / bouwen:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Registratie</w:t></w:r><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> / Authentication</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Login gedeelte</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Home Page</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Product view / overzicht / zoekfunctie</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Order plaatsen</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Admin paneel</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Doel is een shippable product</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Prioriteiten:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Iteration gedeelte werkt</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Zie de backend interactie met bijv. Postman</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Daarop wordt de UI gebouwd</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Met welke feature aan de slag gaan, wat is haalbaar qua timemanagement.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Volgende week Dockerpipeline</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Material Design</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Aan het eind van een shippable product bij iedere iteratie, een backlog over hebben.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Idee: de issues in Trello koppelen aan een commit</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
[void]$target.Range.InsertXML($xml)
